$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Expand the "Parameter" field's instructions text (cell B6) with a note about
# WQX vs. Simple parameter names, matching the updated Meta template guidance.
$newText = "Name of the measured parameter." + [char]10 + "- Note that this can be either the WQX or Simple parameter name.  However, if a parameter is distinguished by Sample Fraction only (i.e. TDP, PON, TDN), then the Simple parameter name must be used here and in all other files (Results, DQO, etc.)"

$ws.Range("B6").Value = $newText
